$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Năm học" (school-year) column so the table gets a second
# year column (2019 / 2020) alongside the existing one, pushing the
# "Họ và tên"/"Lớp"/"Chuyên khoa" columns one slot to the right.
$ws.Columns("B").Insert()

# Re-assert every cell explicitly (header row + the two data rows) so the
# final grid matches exactly, regardless of what the column insert copied.
$ws.Range("A1").Value = "Số hiệu"
$ws.Range("B1").Value = "Năm học"
$ws.Range("C1").Value = "Năm học"
$ws.Range("D1").Value = "Họ và tên"
$ws.Range("E1").Value = "Lớp"
$ws.Range("F1").Value = "Chuyên khoa"

# Row 2 (existing candidate): serial number becomes a text value "49.330"
# formatted as Text (numFmt 49 / "@"); new "2020" year column added.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "49.330"
$ws.Range("B2").Value = 2019
$ws.Range("C2").Value = 2020
$ws.Range("D2").Value = "Nguyen Tuan Nghia"
$ws.Range("E2").Value = "B12D49"
$ws.Range("F2").Value = 3

# Row 3: a new/updated candidate (per the commit message, one more thí
# sinh was added), replacing the previous "Le Dang Quang"/B12D39 entry.
$ws.Range("A3").Value = 49.331000000000003
$ws.Range("B3").Value = 2019
$ws.Range("C3").Value = 2020
$ws.Range("D3").Value = "Nguyen Van Nghia "
$ws.Range("E3").Value = "B12D48"
$ws.Range("F3").Value = 4

# Match the author's final cursor position.
$ws.Range("F3").Select()
